# Actualización automática 2025-06-24 14:30:08
# Updates the sales figures for ALMEIDA CUATIN JHONATHANN CARLOS / PORCELANATO
# across the three report sheets, keeping dependent totals/percentages in sync.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": M18 (PORCELANATO column for this advisor/client row)
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M18").Value = 1324.81

# --- Sheet "VENTA MENSUAL": F18 (row total for this advisor/client) and F31 (grand total)
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F18").Value = 1641.09
$wsMensual.Range("F31").Value = 7454.38

# --- Sheet "CUMPLIMIENTO MENSUAL": PORCELANATO row (16) and TOTAL row (19)
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D16").Value = 4896.65
$wsCumpl.Range("E16").Value = 13901.96
$wsCumpl.Range("F16").Value = 0.2604793652296633

$wsCumpl.Range("D19").Value = 7448.62
$wsCumpl.Range("E19").Value = 22089.17107555787
$wsCumpl.Range("F19").Value = 0.2521725467197726
